$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.377.16"
$ws.Range("E2").Value = "  -1.94%  "
$ws.Range("D3").Value = "1.664.08"
$ws.Range("E3").Value = "  -3.81%  "
$ws.Range("D4").Value = "'0.9989"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.19%  "
$ws.Range("D5").Value = "'235.96"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.57%  "
$ws.Range("D6").Value = "'0.9999"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.06%  "
$ws.Range("D7").Value = "'0.4793"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -4.52%  "
$ws.Range("D8").Value = "'0.2607"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.64%  "
$ws.Range("D9").Value = "'0.06148"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.07%  "
$ws.Range("D10").Value = "'0.07079"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.97%  "
$ws.Range("D11").Value = "1.658.36"
$ws.Range("E11").Value = "  -4.19%  "
$ws.Range("D12").Value = "'14.69"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.70%  "
$ws.Range("D13").Value = "'0.5912"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -8.58%  "
$ws.Range("D14").Value = "'4.380"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -7.46%  "
$ws.Range("D15").Value = "'74.24"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.21%  "
$ws.Range("D16").Value = "'0.9999"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.07%  "
$ws.Range("D17").Value = "'0.9999"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.10%  "
$ws.Range("D18").Value = "25.373.41"
$ws.Range("E18").Value = "  -1.98%  "
$ws.Range("D19").Value = "'0.000006725"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.13%  "
$ws.Range("D20").Value = "'11.41"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.87%  "
$ws.Range("D21").Value = "1.872.11"
$ws.Range("E21").Value = "  -4.25%  "
$ws.Range("D22").Value = "'4.437"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.15%  "
$ws.Range("D23").Value = "'8.651"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.35%  "
$ws.Range("D24").Value = "'5.326"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.57%  "
$ws.Range("D25").Value = "'133.49"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.56%  "
$ws.Range("D26").Value = "'15.07"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.23%  "
$ws.Range("D27").Value = "'1.399"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.16%  "
$ws.Range("D28").Value = "'104.38"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.67%  "
$ws.Range("D29").Value = "'1.692"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.82%  "
$ws.Range("D30").Value = "'3.972"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.51%  "
$ws.Range("D31").Value = "'3.619"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.84%  "
$ws.Range("D32").Value = "'0.07649"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -5.68%  "
$ws.Range("D33").Value = "'0.04392"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -6.71%  "
$ws.Range("D34").Value = "'0.9991"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.06%  "
$ws.Range("D35").Value = "'2.601"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.91%  "
$ws.Range("D36").Value = "'0.6082"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.20%  "
$ws.Range("D37").Value = "'0.9421"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -5.42%  "
$ws.Range("D38").Value = "'2.625"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.19%  "
$ws.Range("D39").Value = "'0.8527"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.75%  "
$ws.Range("D40").Value = "'1.001"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.01%  "
$ws.Range("D41").Value = "'0.01501"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -6.58%  "
$ws.Range("D42").Value = "'1.819"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -6.51%  "
$ws.Range("D43").Value = "'98.48"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.98%  "
$ws.Range("D44").Value = "'0.3761"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.86%  "
$ws.Range("D45").Value = "'4.686"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -6.14%  "
$ws.Range("D46").Value = "'0.1117"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -5.33%  "
$ws.Range("D47").Value = "'6.210"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.86%  "
$ws.Range("D48").Value = "'0.05252"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.52%  "
$ws.Range("D49").Value = "'29.52"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.76%  "
$ws.Range("D50").Value = "'1.211"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.86%  "
$ws.Range("D51").Value = "'1.001"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.15%  "
